$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect text-looking numeric values (e.g. "1.001", "24.931.64") from
# being auto-converted to numbers by forcing a Text number format before
# writing the values.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.907.96"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "1.677.06"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "328.97"
$ws.Range("E5").Value = "  +7.48%  "
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.3658"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("D8").Value = "47.19"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "0.3256"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "1.147"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").Value = "0.07099"
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "6.107"
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").Value = "19.74"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "1.678.20"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "6.636"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "0.00001052"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "0.06607"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").Value = "0.9991"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "79.14"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("D21").Value = "15.98"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "5.939"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "12.85"
$ws.Range("E23").Value = "  +5.61%  "
$ws.Range("D24").Value = "24.947.90"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").Value = "2.455"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "2.421"
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("D27").Value = "148.63"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "18.80"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("D29").Value = "1.863.88"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("D30").Value = "126.12"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "1.187"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").Value = "4.074"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "5.796"
$ws.Range("E33").Value = "  +4.75%  "
$ws.Range("D34").Value = "0.08474"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "1.654"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "5.194"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").Value = "0.02267"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.233"
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06040"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").Value = "0.2099"
$ws.Range("E41").Value = "  +3.32%  "
$ws.Range("D42").Value = "8.276"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "0.5975"
$ws.Range("E44").Value = "  +3.46%  "
$ws.Range("D45").Value = "13.62"
$ws.Range("E45").Value = "  +7.61%  "
$ws.Range("D46").Value = "3.849"
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("D47").Value = "0.5737"
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").Value = "125.59"
$ws.Range("E48").Value = "  +3.41%  "
$ws.Range("D49").Value = "1.966"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "0.07029"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "1.191"
$ws.Range("E51").Value = "  +3.56%  "

# Restore the default (Normal) style so the saved file keeps the same
# cell styling as before the edit (no explicit style index).
$ws.Range("B2:E51").Style = "Normal"

